$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value2 = "ECs"
$ws.Cells.Item(2,2).Value2 = "Thbs2"
$ws.Cells.Item(2,3).Value2 = "Cd47"
$ws.Cells.Item(2,4).Value2 = "ECs"
$ws.Cells.Item(2,5).Value2 = 3
$ws.Cells.Item(2,6).Value2 = 1
$ws.Cells.Item(2,7).Value2 = 1.924015333333333
$ws.Cells.Item(2,8).Value2 = 5.772046
$ws.Cells.Item(2,9).Value2 = 0.009877822204539638
$ws.Cells.Item(2,10).Value2 = 0.009877822204539637
$ws.Cells.Item(2,11).Value2 = 3
$ws.Cells.Item(2,12).Value2 = 1
$ws.Cells.Item(2,13).Value2 = 46.33695966666667
$ws.Cells.Item(2,14).Value2 = 139.010879
$ws.Cells.Item(2,15).Value2 = 0.1993490803952133
$ws.Cells.Item(2,16).Value2 = 0.1993490803952133
$ws.Cells.Item(2,17).Value2 = 89.15302089871489
$ws.Cells.Item(2,18).Value2 = 802.3771880884339
$ws.Cells.Item(2,19).Value2 = 0.001969134772782396
$ws.Cells.Item(2,20).Value2 = 0.001969134772782395

$ws.Cells.Item(3,1).Value2 = "ECs"
$ws.Cells.Item(3,2).Value2 = "Thbs2"
$ws.Cells.Item(3,3).Value2 = "Cd47"
$ws.Cells.Item(3,4).Value2 = "FAPs"
$ws.Cells.Item(3,5).Value2 = 3
$ws.Cells.Item(3,6).Value2 = 1
$ws.Cells.Item(3,7).Value2 = 1.924015333333333
$ws.Cells.Item(3,8).Value2 = 5.772046
$ws.Cells.Item(3,9).Value2 = 0.009877822204539638
$ws.Cells.Item(3,10).Value2 = 0.009877822204539637
$ws.Cells.Item(3,11).Value2 = 3
$ws.Cells.Item(3,12).Value2 = 1
$ws.Cells.Item(3,13).Value2 = 84.50960033333332
$ws.Cells.Item(3,14).Value2 = 253.528801
$ws.Cells.Item(3,15).Value2 = 0.3635739425333109
$ws.Cells.Item(3,16).Value2 = 0.3635739425333109
$ws.Cells.Item(3,17).Value2 = 162.5977668552051
$ws.Cells.Item(3,18).Value2 = 1463.379901696846
$ws.Cells.Item(3,19).Value2 = 0.003591318762547557
$ws.Cells.Item(3,20).Value2 = 0.003591318762547556

$ws.Cells.Item(4,1).Value2 = "ECs"
$ws.Cells.Item(4,2).Value2 = "Thbs2"
$ws.Cells.Item(4,3).Value2 = "Cd47"
$ws.Cells.Item(4,4).Value2 = "M2"
$ws.Cells.Item(4,5).Value2 = 3
$ws.Cells.Item(4,6).Value2 = 1
$ws.Cells.Item(4,7).Value2 = 1.924015333333333
$ws.Cells.Item(4,8).Value2 = 5.772046
$ws.Cells.Item(4,9).Value2 = 0.009877822204539638
$ws.Cells.Item(4,10).Value2 = 0.009877822204539637
$ws.Cells.Item(4,11).Value2 = 3
$ws.Cells.Item(4,12).Value2 = 1
$ws.Cells.Item(4,13).Value2 = 72.52790466666666
$ws.Cells.Item(4,14).Value2 = 217.583714
$ws.Cells.Item(4,15).Value2 = 0.3120267536390091
$ws.Cells.Item(4,16).Value2 = 0.3120267536390091
$ws.Cells.Item(4,17).Value2 = 139.5448006732049
$ws.Cells.Item(4,18).Value2 = 1255.903206058844
$ws.Cells.Item(4,19).Value2 = 0.003082144795505823
$ws.Cells.Item(4,20).Value2 = 0.003082144795505822

$ws.Cells.Item(5,1).Value2 = "ECs"
$ws.Cells.Item(5,2).Value2 = "Thbs2"
$ws.Cells.Item(5,3).Value2 = "Cd47"
$ws.Cells.Item(5,4).Value2 = "sCs"
$ws.Cells.Item(5,5).Value2 = 3
$ws.Cells.Item(5,6).Value2 = 1
$ws.Cells.Item(5,7).Value2 = 1.924015333333333
$ws.Cells.Item(5,8).Value2 = 5.772046
$ws.Cells.Item(5,9).Value2 = 0.009877822204539638
$ws.Cells.Item(5,10).Value2 = 0.009877822204539637
$ws.Cells.Item(5,11).Value2 = 3
$ws.Cells.Item(5,12).Value2 = 1
$ws.Cells.Item(5,13).Value2 = 29.06683666666666
$ws.Cells.Item(5,14).Value2 = 87.20050999999998
$ws.Cells.Item(5,15).Value2 = 0.1250502234324667
$ws.Cells.Item(5,16).Value2 = 0.1250502234324667
$ws.Cells.Item(5,17).Value2 = 55.92503943816221
$ws.Cells.Item(5,18).Value2 = 503.3253549434598
$ws.Cells.Item(5,19).Value2 = 0.001235223873703863
$ws.Cells.Item(5,20).Value2 = 0.001235223873703863

$ws.Cells.Item(6,1).Value2 = "FAPs"
$ws.Cells.Item(6,2).Value2 = "Thbs2"
$ws.Cells.Item(6,3).Value2 = "Cd47"
$ws.Cells.Item(6,4).Value2 = "ECs"
$ws.Cells.Item(6,5).Value2 = 3
$ws.Cells.Item(6,6).Value2 = 1
$ws.Cells.Item(6,7).Value2 = 188.2309416666667
$ws.Cells.Item(6,8).Value2 = 564.692825
$ws.Cells.Item(6,9).Value2 = 0.9663705600283187
$ws.Cells.Item(6,10).Value2 = 0.9663705600283184
$ws.Cells.Item(6,11).Value2 = 3
$ws.Cells.Item(6,12).Value2 = 1
$ws.Cells.Item(6,13).Value2 = 46.33695966666667
$ws.Cells.Item(6,14).Value2 = 139.010879
$ws.Cells.Item(6,15).Value2 = 0.1993490803952133
$ws.Cells.Item(6,16).Value2 = 0.1993490803952133
$ws.Cells.Item(6,17).Value2 = 8722.04955202702
$ws.Cells.Item(6,18).Value2 = 78498.44596824316
$ws.Cells.Item(6,19).Value2 = 0.1926450824626526
$ws.Cells.Item(6,20).Value2 = 0.1926450824626526

$ws.Cells.Item(7,1).Value2 = "FAPs"
$ws.Cells.Item(7,2).Value2 = "Thbs2"
$ws.Cells.Item(7,3).Value2 = "Cd47"
$ws.Cells.Item(7,4).Value2 = "FAPs"
$ws.Cells.Item(7,5).Value2 = 3
$ws.Cells.Item(7,6).Value2 = 1
$ws.Cells.Item(7,7).Value2 = 188.2309416666667
$ws.Cells.Item(7,8).Value2 = 564.692825
$ws.Cells.Item(7,9).Value2 = 0.9663705600283187
$ws.Cells.Item(7,10).Value2 = 0.9663705600283184
$ws.Cells.Item(7,11).Value2 = 3
$ws.Cells.Item(7,12).Value2 = 1
$ws.Cells.Item(7,13).Value2 = 84.50960033333332
$ws.Cells.Item(7,14).Value2 = 253.528801
$ws.Cells.Item(7,15).Value2 = 0.3635739425333109
$ws.Cells.Item(7,16).Value2 = 0.3635739425333109
$ws.Cells.Item(7,17).Value2 = 15907.32165061698
$ws.Cells.Item(7,18).Value2 = 143165.8948555528
$ws.Cells.Item(7,19).Value2 = 0.3513471544576194
$ws.Cells.Item(7,20).Value2 = 0.3513471544576193

$ws.Cells.Item(8,1).Value2 = "FAPs"
$ws.Cells.Item(8,2).Value2 = "Thbs2"
$ws.Cells.Item(8,3).Value2 = "Cd47"
$ws.Cells.Item(8,4).Value2 = "M2"
$ws.Cells.Item(8,5).Value2 = 3
$ws.Cells.Item(8,6).Value2 = 1
$ws.Cells.Item(8,7).Value2 = 188.2309416666667
$ws.Cells.Item(8,8).Value2 = 564.692825
$ws.Cells.Item(8,9).Value2 = 0.9663705600283187
$ws.Cells.Item(8,10).Value2 = 0.9663705600283184
$ws.Cells.Item(8,11).Value2 = 3
$ws.Cells.Item(8,12).Value2 = 1
$ws.Cells.Item(8,13).Value2 = 72.52790466666666
$ws.Cells.Item(8,14).Value2 = 217.583714
$ws.Cells.Item(8,15).Value2 = 0.3120267536390091
$ws.Cells.Item(8,16).Value2 = 0.3120267536390091
$ws.Cells.Item(8,17).Value2 = 13651.99579251689
$ws.Cells.Item(8,18).Value2 = 122867.962132652
$ws.Cells.Item(8,19).Value2 = 0.3015334686579474
$ws.Cells.Item(8,20).Value2 = 0.3015334686579473

$ws.Cells.Item(9,1).Value2 = "FAPs"
$ws.Cells.Item(9,2).Value2 = "Thbs2"
$ws.Cells.Item(9,3).Value2 = "Cd47"
$ws.Cells.Item(9,4).Value2 = "sCs"
$ws.Cells.Item(9,5).Value2 = 3
$ws.Cells.Item(9,6).Value2 = 1
$ws.Cells.Item(9,7).Value2 = 188.2309416666667
$ws.Cells.Item(9,8).Value2 = 564.692825
$ws.Cells.Item(9,9).Value2 = 0.9663705600283187
$ws.Cells.Item(9,10).Value2 = 0.9663705600283184
$ws.Cells.Item(9,11).Value2 = 3
$ws.Cells.Item(9,12).Value2 = 1
$ws.Cells.Item(9,13).Value2 = 29.06683666666666
$ws.Cells.Item(9,14).Value2 = 87.20050999999998
$ws.Cells.Item(9,15).Value2 = 0.1250502234324667
$ws.Cells.Item(9,16).Value2 = 0.1250502234324667
$ws.Cells.Item(9,17).Value2 = 5471.27803703786
$ws.Cells.Item(9,18).Value2 = 49241.50233334074
$ws.Cells.Item(9,19).Value2 = 0.1208448544500993
$ws.Cells.Item(9,20).Value2 = 0.1208448544500992

$ws.Cells.Item(10,1).Value2 = "M2"
$ws.Cells.Item(10,2).Value2 = "Thbs2"
$ws.Cells.Item(10,3).Value2 = "Cd47"
$ws.Cells.Item(10,4).Value2 = "ECs"
$ws.Cells.Item(10,5).Value2 = 1
$ws.Cells.Item(10,6).Value2 = 0.3333333333333333
$ws.Cells.Item(10,7).Value2 = 0.08498099999999999
$ws.Cells.Item(10,8).Value2 = 0.254943
$ws.Cells.Item(10,9).Value2 = 0.000436289251037145
$ws.Cells.Item(10,10).Value2 = 0.000436289251037145
$ws.Cells.Item(10,11).Value2 = 3
$ws.Cells.Item(10,12).Value2 = 1
$ws.Cells.Item(10,13).Value2 = 46.33695966666667
$ws.Cells.Item(10,14).Value2 = 139.010879
$ws.Cells.Item(10,15).Value2 = 0.1993490803952133
$ws.Cells.Item(10,16).Value2 = 0.1993490803952133
$ws.Cells.Item(10,17).Value2 = 3.937761169432999
$ws.Cells.Item(10,18).Value2 = 35.439850524897
$ws.Cells.Item(10,19).Value2 = 0.00008697386098057123
$ws.Cells.Item(10,20).Value2 = 0.0000869738609805712

$ws.Cells.Item(11,1).Value2 = "M2"
$ws.Cells.Item(11,2).Value2 = "Thbs2"
$ws.Cells.Item(11,3).Value2 = "Cd47"
$ws.Cells.Item(11,4).Value2 = "FAPs"
$ws.Cells.Item(11,5).Value2 = 1
$ws.Cells.Item(11,6).Value2 = 0.3333333333333333
$ws.Cells.Item(11,7).Value2 = 0.08498099999999999
$ws.Cells.Item(11,8).Value2 = 0.254943
$ws.Cells.Item(11,9).Value2 = 0.000436289251037145
$ws.Cells.Item(11,10).Value2 = 0.000436289251037145
$ws.Cells.Item(11,11).Value2 = 3
$ws.Cells.Item(11,12).Value2 = 1
$ws.Cells.Item(11,13).Value2 = 84.50960033333332
$ws.Cells.Item(11,14).Value2 = 253.528801
$ws.Cells.Item(11,15).Value2 = 0.3635739425333109
$ws.Cells.Item(11,16).Value2 = 0.3635739425333109
$ws.Cells.Item(11,17).Value2 = 7.181710345926998
$ws.Cells.Item(11,18).Value2 = 64.63539311334299
$ws.Cells.Item(11,19).Value2 = 0.0001586234030844802
$ws.Cells.Item(11,20).Value2 = 0.0001586234030844802

$ws.Cells.Item(12,1).Value2 = "M2"
$ws.Cells.Item(12,2).Value2 = "Thbs2"
$ws.Cells.Item(12,3).Value2 = "Cd47"
$ws.Cells.Item(12,4).Value2 = "M2"
$ws.Cells.Item(12,5).Value2 = 1
$ws.Cells.Item(12,6).Value2 = 0.3333333333333333
$ws.Cells.Item(12,7).Value2 = 0.08498099999999999
$ws.Cells.Item(12,8).Value2 = 0.254943
$ws.Cells.Item(12,9).Value2 = 0.000436289251037145
$ws.Cells.Item(12,10).Value2 = 0.000436289251037145
$ws.Cells.Item(12,11).Value2 = 3
$ws.Cells.Item(12,12).Value2 = 1
$ws.Cells.Item(12,13).Value2 = 72.52790466666666
$ws.Cells.Item(12,14).Value2 = 217.583714
$ws.Cells.Item(12,15).Value2 = 0.3120267536390091
$ws.Cells.Item(12,16).Value2 = 0.3120267536390091
$ws.Cells.Item(12,17).Value2 = 6.163493866477999
$ws.Cells.Item(12,18).Value2 = 55.47144479830199
$ws.Cells.Item(12,19).Value2 = 0.000136133918648715
$ws.Cells.Item(12,20).Value2 = 0.000136133918648715

$ws.Cells.Item(13,1).Value2 = "M2"
$ws.Cells.Item(13,2).Value2 = "Thbs2"
$ws.Cells.Item(13,3).Value2 = "Cd47"
$ws.Cells.Item(13,4).Value2 = "sCs"
$ws.Cells.Item(13,5).Value2 = 1
$ws.Cells.Item(13,6).Value2 = 0.3333333333333333
$ws.Cells.Item(13,7).Value2 = 0.08498099999999999
$ws.Cells.Item(13,8).Value2 = 0.254943
$ws.Cells.Item(13,9).Value2 = 0.000436289251037145
$ws.Cells.Item(13,10).Value2 = 0.000436289251037145
$ws.Cells.Item(13,11).Value2 = 3
$ws.Cells.Item(13,12).Value2 = 1
$ws.Cells.Item(13,13).Value2 = 29.06683666666666
$ws.Cells.Item(13,14).Value2 = 87.20050999999998
$ws.Cells.Item(13,15).Value2 = 0.1250502234324667
$ws.Cells.Item(13,16).Value2 = 0.1250502234324667
$ws.Cells.Item(13,17).Value2 = 2.470128846769999
$ws.Cells.Item(13,18).Value2 = 22.23115962092999
$ws.Cells.Item(13,19).Value2 = 0.00005455806832337856
$ws.Cells.Item(13,20).Value2 = 0.00005455806832337855

$ws.Cells.Item(14,1).Value2 = "sCs"
$ws.Cells.Item(14,2).Value2 = "Thbs2"
$ws.Cells.Item(14,3).Value2 = "Cd47"
$ws.Cells.Item(14,4).Value2 = "ECs"
$ws.Cells.Item(14,5).Value2 = 3
$ws.Cells.Item(14,6).Value2 = 1
$ws.Cells.Item(14,7).Value2 = 4.541390666666667
$ws.Cells.Item(14,8).Value2 = 13.624172
$ws.Cells.Item(14,9).Value2 = 0.02331532851610455
$ws.Cells.Item(14,10).Value2 = 0.02331532851610455
$ws.Cells.Item(14,11).Value2 = 3
$ws.Cells.Item(14,12).Value2 = 1
$ws.Cells.Item(14,13).Value2 = 46.33695966666667
$ws.Cells.Item(14,14).Value2 = 139.010879
$ws.Cells.Item(14,15).Value2 = 0.1993490803952133
$ws.Cells.Item(14,16).Value2 = 0.1993490803952133
$ws.Cells.Item(14,17).Value2 = 210.4342361519098
$ws.Cells.Item(14,18).Value2 = 1893.908125367188
$ws.Cells.Item(14,19).Value2 = 0.004647889298797736
$ws.Cells.Item(14,20).Value2 = 0.004647889298797735

$ws.Cells.Item(15,1).Value2 = "sCs"
$ws.Cells.Item(15,2).Value2 = "Thbs2"
$ws.Cells.Item(15,3).Value2 = "Cd47"
$ws.Cells.Item(15,4).Value2 = "FAPs"
$ws.Cells.Item(15,5).Value2 = 3
$ws.Cells.Item(15,6).Value2 = 1
$ws.Cells.Item(15,7).Value2 = 4.541390666666667
$ws.Cells.Item(15,8).Value2 = 13.624172
$ws.Cells.Item(15,9).Value2 = 0.02331532851610455
$ws.Cells.Item(15,10).Value2 = 0.02331532851610455
$ws.Cells.Item(15,11).Value2 = 3
$ws.Cells.Item(15,12).Value2 = 1
$ws.Cells.Item(15,13).Value2 = 84.50960033333332
$ws.Cells.Item(15,14).Value2 = 253.528801
$ws.Cells.Item(15,15).Value2 = 0.3635739425333109
$ws.Cells.Item(15,16).Value2 = 0.3635739425333109
$ws.Cells.Item(15,17).Value2 = 383.7911101975302
$ws.Cells.Item(15,18).Value2 = 3454.119991777772
$ws.Cells.Item(15,19).Value2 = 0.008476845910059462
$ws.Cells.Item(15,20).Value2 = 0.00847684591005946

$ws.Cells.Item(16,1).Value2 = "sCs"
$ws.Cells.Item(16,2).Value2 = "Thbs2"
$ws.Cells.Item(16,3).Value2 = "Cd47"
$ws.Cells.Item(16,4).Value2 = "M2"
$ws.Cells.Item(16,5).Value2 = 3
$ws.Cells.Item(16,6).Value2 = 1
$ws.Cells.Item(16,7).Value2 = 4.541390666666667
$ws.Cells.Item(16,8).Value2 = 13.624172
$ws.Cells.Item(16,9).Value2 = 0.02331532851610455
$ws.Cells.Item(16,10).Value2 = 0.02331532851610455
$ws.Cells.Item(16,11).Value2 = 3
$ws.Cells.Item(16,12).Value2 = 1
$ws.Cells.Item(16,13).Value2 = 72.52790466666666
$ws.Cells.Item(16,14).Value2 = 217.583714
$ws.Cells.Item(16,15).Value2 = 0.3120267536390091
$ws.Cells.Item(16,16).Value2 = 0.3120267536390091
$ws.Cells.Item(16,17).Value2 = 329.3775493260897
$ws.Cells.Item(16,18).Value2 = 2964.397943934808
$ws.Cells.Item(16,19).Value2 = 0.007275006266907118
$ws.Cells.Item(16,20).Value2 = 0.007275006266907118

$ws.Cells.Item(17,1).Value2 = "sCs"
$ws.Cells.Item(17,2).Value2 = "Thbs2"
$ws.Cells.Item(17,3).Value2 = "Cd47"
$ws.Cells.Item(17,4).Value2 = "sCs"
$ws.Cells.Item(17,5).Value2 = 3
$ws.Cells.Item(17,6).Value2 = 1
$ws.Cells.Item(17,7).Value2 = 4.541390666666667
$ws.Cells.Item(17,8).Value2 = 13.624172
$ws.Cells.Item(17,9).Value2 = 0.02331532851610455
$ws.Cells.Item(17,10).Value2 = 0.02331532851610455
$ws.Cells.Item(17,11).Value2 = 3
$ws.Cells.Item(17,12).Value2 = 1
$ws.Cells.Item(17,13).Value2 = 29.06683666666666
$ws.Cells.Item(17,14).Value2 = 87.20050999999998
$ws.Cells.Item(17,15).Value2 = 0.1250502234324667
$ws.Cells.Item(17,16).Value2 = 0.1250502234324667
$ws.Cells.Item(17,17).Value2 = 132.0038607475244
$ws.Cells.Item(17,18).Value2 = 1188.03474672772
$ws.Cells.Item(17,19).Value2 = 0.002915587040340238
$ws.Cells.Item(17,20).Value2 = 0.002915587040340238

